$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (2007年 data) and shift the remaining rows (2010, 2012, 2015) up.
$ws.Rows.Item(2).Delete()
